# Percentage_of_Instruction.xlsx - "Working more on percent instruction"
#
# 1. Sheet1 ("Sheet1" -> "VOCALS"): add a new row 16 ("&gt;= 0 dB") of data.
# 2. Add a new sheet "VOCALS x 2" (copy of the VOCALS layout, with new data
#    in the 1/0.1, 8/0.1, 10/0.1 rows, and an extra blank ">= 0 dB" row).
# 3. Add a new sheet "VOCALS x 10" (same layout again, with a new leading
#    "0/0.1" row, new data, and an extra blank ">= 0 dB" row).

$wb = $excel.ActiveWorkbook
$vocals = $wb.Worksheets.Item(1)
$vocals.Name = "VOCALS"

# ---------------------------------------------------------------------------
# 1. VOCALS: append row 16
# ---------------------------------------------------------------------------
$vocals.Range("A16").Value = ">= 0 dB"

$vocals.Range("B16").Value = 0.4785
$vocals.Range("B16").NumberFormat = $vocals.Range("B5").NumberFormat

$vocals.Range("C16").Value = 0.5079
$vocals.Range("C16").NumberFormat = $vocals.Range("C5").NumberFormat

$vocals.Range("D16").Value = 0.5023
$vocals.Range("D16").NumberFormat = $vocals.Range("D5").NumberFormat

$vocals.Range("B16").Select()

# ---------------------------------------------------------------------------
# 2. New sheet: "VOCALS x 2"
# ---------------------------------------------------------------------------
$vocals2 = $wb.Worksheets.Add($null, $vocals)
$vocals2.Name = "VOCALS x 2"
$vocals2.Columns.Item(1).ColumnWidth = 20.7109375

$vocals2.Range("A1").Value = "Using video_#_vocals_x2.wav files after ran through spleeter and decibels * 2 in ffmpeg"
$vocals2.Range("A1").WrapText = $true
$vocals2.Rows.Item(1).RowHeight = 42

$vocals2.Range("A3").Value = "ACTUAL"

$vocals2.Range("A4").Value = "LABELS"
$vocals2.Range("B4").Value = "VIDEO_1"
$vocals2.Range("C4").Value = "VIDEO_2"
$vocals2.Range("D4").Value = "VIDEO_3"

$vocals2.Range("A5").Value = "1/0.1"
$vocals2.Range("B5").Value = 0.4999
$vocals2.Range("B5").NumberFormat = "0.00%"
$vocals2.Range("C5").NumberFormat = "0.00%"
$vocals2.Range("D5").NumberFormat = "0.00%"

$vocals2.Range("A6").Value = "8/0.1"
$vocals2.Range("B6").Value = 1.1402
$vocals2.Range("B6").NumberFormat = "0.00%"
$vocals2.Range("C6").NumberFormat = "0.00%"
$vocals2.Range("D6").NumberFormat = "0.00%"

$vocals2.Range("A7").Value = "10/0.1"
$vocals2.Range("B7").Value = 1.8193
$vocals2.Range("B7").NumberFormat = "0.00%"
$vocals2.Range("C7").NumberFormat = "0.00%"
$vocals2.Range("D7").NumberFormat = "0.00%"

$vocals2.Range("A10").Value = "Percent_of_Instruction"

$vocals2.Range("A11").Value = "LABELS"
$vocals2.Range("B11").Value = "VIDEO_1"
$vocals2.Range("C11").Value = "VIDEO_2"
$vocals2.Range("D11").Value = "VIDEO_3"

$vocals2.Range("A12").Value = "20 dB - 60 dB"
$vocals2.Range("B12").NumberFormat = "0.00%"
$vocals2.Range("C12").NumberFormat = "0.00%"
$vocals2.Range("D12").NumberFormat = "0.00%"

$vocals2.Range("A13").Value = "0 dB - 60 dB"
$vocals2.Range("B13").NumberFormat = "0.00%"
$vocals2.Range("C13").NumberFormat = "0.00%"
$vocals2.Range("D13").NumberFormat = "0.00%"

$vocals2.Range("A14").Value = "20 dB - 80 dB"
$vocals2.Range("B14").NumberFormat = "0.00%"
$vocals2.Range("C14").NumberFormat = "0.00%"
$vocals2.Range("D14").NumberFormat = "0.00%"

$vocals2.Range("A15").Value = "0 dB - 80 dB"
$vocals2.Range("B15").NumberFormat = "0.00%"
$vocals2.Range("C15").NumberFormat = "0.00%"
$vocals2.Range("D15").NumberFormat = "0.00%"

$vocals2.Range("A16").Value = ">= 0 dB"
$vocals2.Range("B16").NumberFormat = "0.00%"
$vocals2.Range("C16").NumberFormat = "0.00%"
$vocals2.Range("D16").NumberFormat = "0.00%"

$vocals2.Range("B15").Select()
$excel.ActiveWindow.Zoom = 205
$vocals2.Activate()

# ---------------------------------------------------------------------------
# 3. New sheet: "VOCALS x 10"
# ---------------------------------------------------------------------------
$vocals10 = $wb.Worksheets.Add($null, $vocals2)
$vocals10.Name = "VOCALS x 10"
$vocals10.Columns.Item(1).ColumnWidth = 26.85546875

$vocals10.Range("A1").Value = "Using video_#_vocals_x10.wav files after ran through spleeter and decibels * 2 in ffmpeg"
$vocals10.Range("A1").WrapText = $true
$vocals10.Rows.Item(1).RowHeight = 42

$vocals10.Range("A3").Value = "ACTUAL"

$vocals10.Range("A4").Value = "LABELS"
$vocals10.Range("B4").Value = "VIDEO_1"
$vocals10.Range("C4").Value = "VIDEO_2"
$vocals10.Range("D4").Value = "VIDEO_3"

$vocals10.Range("A5").Value = "0/0.1"
$vocals10.Range("B5").Value = 0.749
$vocals10.Range("B5").NumberFormat = "0.00%"

$vocals10.Range("A6").Value = "1/0.1"
$vocals10.Range("B6").Value = 1.5166
$vocals10.Range("B6").NumberFormat = "0.00%"
$vocals10.Range("C6").NumberFormat = "0.00%"
$vocals10.Range("D6").NumberFormat = "0.00%"

$vocals10.Range("A7").Value = "8/0.1"
$vocals10.Range("B7").NumberFormat = "0.00%"
$vocals10.Range("C7").NumberFormat = "0.00%"
$vocals10.Range("D7").NumberFormat = "0.00%"

$vocals10.Range("A8").Value = "10/0.1"
$vocals10.Range("B8").NumberFormat = "0.00%"
$vocals10.Range("C8").NumberFormat = "0.00%"
$vocals10.Range("D8").NumberFormat = "0.00%"

$vocals10.Range("A11").Value = "Percent_of_Instruction"

$vocals10.Range("A12").Value = "LABELS"
$vocals10.Range("B12").Value = "VIDEO_1"
$vocals10.Range("C12").Value = "VIDEO_2"
$vocals10.Range("D12").Value = "VIDEO_3"

$vocals10.Range("A13").Value = "20 dB - 60 dB"
$vocals10.Range("B13").NumberFormat = "0.00%"
$vocals10.Range("C13").NumberFormat = "0.00%"
$vocals10.Range("D13").NumberFormat = "0.00%"

$vocals10.Range("A14").Value = "0 dB - 60 dB"
$vocals10.Range("B14").NumberFormat = "0.00%"
$vocals10.Range("C14").NumberFormat = "0.00%"
$vocals10.Range("D14").NumberFormat = "0.00%"

$vocals10.Range("A15").Value = "20 dB - 80 dB"
$vocals10.Range("B15").NumberFormat = "0.00%"
$vocals10.Range("C15").NumberFormat = "0.00%"
$vocals10.Range("D15").NumberFormat = "0.00%"

$vocals10.Range("A16").Value = "0 dB - 80 dB"
$vocals10.Range("B16").NumberFormat = "0.00%"
$vocals10.Range("C16").NumberFormat = "0.00%"
$vocals10.Range("D16").NumberFormat = "0.00%"

$vocals10.Range("A17").Value = ">= 0 dB"
$vocals10.Range("B17").NumberFormat = "0.00%"
$vocals10.Range("C17").NumberFormat = "0.00%"
$vocals10.Range("D17").NumberFormat = "0.00%"

$vocals10.Range("C8").Select()

# Re-activate "VOCALS x 2" as the tab that was selected when saved.
$vocals2.Activate()
